$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52
$prevRow = 51

# Copy formatting from the previous row's date cell so the new row matches
# the existing style (numeric date format) used throughout column A.
$ws.Cells.Item($prevRow, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($row, 1).Value = 45956
$ws.Cells.Item($row, 2).Value = "21,7048"
$ws.Cells.Item($row, 3).Value = "15,5758"
$ws.Cells.Item($row, 4).Value = "15,5156"
$ws.Cells.Item($row, 5).Value = "15,5156"
